$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.111.49"
$ws.Range("E2").Value = "  +0.36%  "
# Row 3
$ws.Range("D3").Value = "1.851.63"
$ws.Range("E3").Value = "  +1.86%  "
# Row 4
$ws.Range("E4").Value = "  +0.60%  "
# Row 5
$ws.Range("D5").Value = "'237.78"
$ws.Range("E5").Value = "  +3.37%  "
# Row 6
$ws.Range("D6").Value = "'0.620"
$ws.Range("E6").Value = "  +0.94%  "
# Row 7
$ws.Range("E7").Value = "  +0.50%  "
# Row 8
$ws.Range("D8").Value = "'41.90"
$ws.Range("E8").Value = "  +4.53%  "
# Row 9
$ws.Range("D9").Value = "'0.326"
$ws.Range("E9").Value = "  +0.93%  "
# Row 10
$ws.Range("D10").Value = "'0.0691"
$ws.Range("E10").Value = "  +1.22%  "
# Row 11
$ws.Range("E11").Value = "  -0.10%  "
# Row 12
$ws.Range("D12").Value = "2.120.30"
$ws.Range("E12").Value = "  +1.92%  "
# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'11.36"
$ws.Range("E13").Value = "  +0.58%  "
# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.835.68"
$ws.Range("E14").Value = "  +0.96%  "
# Row 15
$ws.Range("D15").Value = "'0.672"
$ws.Range("E15").Value = "  +0.97%  "
# Row 16
$ws.Range("D16").Value = "'4.72"
$ws.Range("E16").Value = "  +2.44%  "
# Row 17
$ws.Range("D17").Value = "35.084.14"
$ws.Range("E17").Value = "  +0.29%  "
# Row 18
$ws.Range("D18").Value = "'69.94"
$ws.Range("E18").Value = "  +0.54%  "
# Row 19
$ws.Range("D19").Value = "0.0₃0789"
$ws.Range("E19").Value = "  +0.58%  "
# Row 20
$ws.Range("D20").Value = "'240.10"
$ws.Range("E20").Value = "  -0.22%  "
# Row 21
$ws.Range("D21").Value = "'12.13"
$ws.Range("E21").Value = "  +0.41%  "
# Row 22
$ws.Range("D22").Value = "'4.71"
$ws.Range("E22").Value = "  +1.21%  "
# Row 23
$ws.Range("E23").Value = "  +0.45%  "
# Row 24
$ws.Range("E24").Value = "  +0.22%  "
# Row 25
$ws.Range("D25").Value = "'169.22"
$ws.Range("E25").Value = "  -2.35%  "
# Row 26
$ws.Range("D26").Value = "'7.99"
$ws.Range("E26").Value = "  +2.29%  "
# Row 27
$ws.Range("D27").Value = "'1.84"
$ws.Range("E27").Value = "  +21.49%  "
# Row 28
$ws.Range("D28").Value = "'17.55"
$ws.Range("E28").Value = "  +1.25%  "
# Row 29
$ws.Range("E29").Value = "  +0.39%  "
# Row 30
$ws.Range("E30").Value = "  +0.54%  "
# Row 31
$ws.Range("D31").Value = "'0.0551"
$ws.Range("E31").Value = "  +0.36%  "
# Row 32
$ws.Range("D32").Value = "'3.97"
$ws.Range("E32").Value = "  -0.39%  "
# Row 33
$ws.Range("D33").Value = "'3.99"
$ws.Range("E33").Value = "  +0.88%  "
# Row 34
$ws.Range("D34").Value = "'1.73"
$ws.Range("E34").Value = "  +27.38%  "
# Row 35
$ws.Range("E35").Value = "  +9.31%  "
# Row 36
$ws.Range("E36").Value = "  +15.08%  "
# Row 37
$ws.Range("E37").Value = "  +3.72%  "
# Row 38
$ws.Range("E38").Value = "  +9.09%  "
# Row 39
$ws.Range("E39").Value = "  +3.87%  "
# Row 40
$ws.Range("D40").Value = "'89.76"
$ws.Range("E40").Value = "  -2.90%  "
# Row 41
$ws.Range("D41").Value = "1.341.26"
$ws.Range("E41").Value = "  +0.17%  "
# Row 42
$ws.Range("D42").Value = "'13.30"
$ws.Range("E42").Value = "  +54.51%  "
# Row 43
$ws.Range("D43").Value = "'14.75"
$ws.Range("E43").Value = "  +1.05%  "
# Row 44
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "'2.45"
$ws.Range("E44").Value = "  +1.22%  "
# Row 45
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'2.29"
$ws.Range("E45").Value = "  +1.43%  "
# Row 46
$ws.Range("D46").Value = "'0.0555"
$ws.Range("E46").Value = "  +6.55%  "
# Row 47
$ws.Range("E47").Value = "  -0.49%  "
# Row 48
$ws.Range("D48").Value = "'6.45"
$ws.Range("E48").Value = "  +3.92%  "
# Row 49
$ws.Range("D49").Value = "2.030.86"
$ws.Range("E49").Value = "  +1.68%  "
# Row 50
$ws.Range("E50").Value = "  +1.60%  "
# Row 51
$ws.Range("E51").Value = "  +0.49%  "
